# BUG: Don't extract header names if none specified (gh-11733)
#
# Add a new worksheet "index_col_none" at the end of the workbook. It holds a
# two-row MultiIndex-style column header (top level "A"/"A"/"B"/"B", second
# level "key"/"val"/"key"/"val") followed by two identical data rows
# (1, 2, 3, 4), used by pandas to verify that read_excel(header=[0, 1])
# does not try to parse a header "names" row when none was supplied.

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing worksheet, and becomes the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "index_col_none"

# Top header row: A, A, B, B
$newSheet.Range("A1").Value = "A"
$newSheet.Range("B1").Value = "A"
$newSheet.Range("C1").Value = "B"
$newSheet.Range("D1").Value = "B"

# Second header row: key, val, key, val
$newSheet.Range("A2").Value = "key"
$newSheet.Range("B2").Value = "val"
$newSheet.Range("C2").Value = "key"
$newSheet.Range("D2").Value = "val"

# Two identical data rows: 1, 2, 3, 4
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = 2
$newSheet.Range("C3").Value = 3
$newSheet.Range("D3").Value = 4

$newSheet.Range("A4").Value = 1
$newSheet.Range("B4").Value = 2
$newSheet.Range("C4").Value = 3
$newSheet.Range("D4").Value = 4

# Center every cell, bold the two header rows (matches the style used for
# the other "key"/"val"-style headers already in this workbook).
$newSheet.Range("A1:D4").HorizontalAlignment = -4108
$newSheet.Range("A1:D2").Font.Bold = $true
